# Insert a new "Title and Content" slide at position 4 (right after the
# "Keyword" slide, before "Cac cong cu tien ich"), reviewing Tran Kim Trung's
# week 1 HTML/Flask-form topics. All later slides simply shift down by one.

$p = $ppt.ActivePresentation

# ppLayoutText = 2 -> same "Title and Content" custom layout used by the
# other content slides in this deck.
$newSlide = $p.Slides.Add(4, 2)

# Title placeholder
$newSlide.Shapes.Item(1).TextFrame.TextRange.Text = "Tran Kim Trung"

# Body / content placeholder
$body = $newSlide.Shapes.Item(2).TextFrame.TextRange
$body.Text = "HTML`rhttps://www.w3schools.com/`rCách gửi dữ liệu từ form -> server`rHow to get data from html form flask`rYoutube"

$body.Paragraphs(2,1).IndentLevel = 2
$body.Paragraphs(2,1).ActionSettings.Item(1).Hyperlink.Address = "https://www.w3schools.com/"
$body.Paragraphs(4,1).IndentLevel = 2
$body.Paragraphs(5,1).IndentLevel = 3
